# Added computation times comparison.
# Fill in the missing "time_numerical_min" values for the n_runs=160 and
# n_runs=320 rows, then leave the selection on B9 (matching the author's
# final cursor position after typing the new data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 15.8
$ws.Range("B8").Value = 29.3

[void]$ws.Range("B9").Select()
